# ============================================================================
# Update "想去人数" (want-to-go count, column F) figures across the four
# sheets, drop the "凹凸世界八周年" row from 本地生活 (local life), and add
# the new "春日计划2024" event into 全部类型 (all types) ahead of the
# existing "世界计划25时" row (pushing everything below it down by one row).
# ============================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (Exhibitions) - column F value bumps only, rows unchanged.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$exhibitionUpdates = @(
    @(5, 1003),
    @(7, 2515),
    @(9, 1227),
    @(10, 902),
    @(11, 602),
    @(12, 902),
    @(13, 1116),
    @(17, 721),
    @(18, 763),
    @(19, 185),
    @(20, 479),
    @(21, 1107),
    @(22, 83),
    @(23, 588),
    @(25, 216),
    @(26, 297),
    @(29, 290),
    @(30, 4234),
    @(31, 479),
    @(36, 140),
    @(37, 1592),
    @(38, 437),
    @(40, 82),
    @(41, 137),
    @(42, 67),
    @(44, 124),
    @(45, 130),
    @(46, 93),
    @(48, 87)
)
foreach ($pair in $exhibitionUpdates) {
    $ws1.Cells.Item($pair[0], 6).Value = $pair[1]
}

# ---------------------------------------------------------------------------
# Sheet "演出" (Performances) - column F value bumps only, rows unchanged.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$performanceUpdates = @(
    @(5, 7),
    @(8, 10),
    @(11, 182),
    @(13, 12),
    @(15, 25),
    @(16, 181)
)
foreach ($pair in $performanceUpdates) {
    $ws2.Cells.Item($pair[0], 6).Value = $pair[1]
}

# ---------------------------------------------------------------------------
# Sheet "本地生活" (Local life) - bump column F, then drop row 4
# (广州·凹凸世界八周年 夏日特调主题嘉年华), shrinking the sheet to 3 rows.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$localLifeUpdates = @(
    @(2, 2257),
    @(3, 723)
)
foreach ($pair in $localLifeUpdates) {
    $ws3.Cells.Item($pair[0], 6).Value = $pair[1]
}
$ws3.Rows.Item(4).Delete()

# ---------------------------------------------------------------------------
# Sheet "全部类型" (All types) - bump column F for the rows that keep their
# position (2-40), then insert a new row at 41 for "广州·春日计划2024——特别
# 二次元不插电音乐会", which pushes the old rows 41-49 down to 42-50.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$allTypesUpdates = @(
    @(2, 2257),
    @(3, 723),
    @(7, 1003),
    @(8, 2515),
    @(10, 1227),
    @(11, 902),
    @(12, 602),
    @(13, 902),
    @(14, 1116),
    @(18, 721),
    @(19, 7),
    @(21, 763),
    @(22, 185),
    @(23, 479),
    @(24, 1107),
    @(26, 83),
    @(27, 588),
    @(29, 216),
    @(30, 297),
    @(33, 4234),
    @(34, 182),
    @(35, 479),
    @(38, 140),
    @(39, 1592),
    @(40, 437)
)
foreach ($pair in $allTypesUpdates) {
    $ws4.Cells.Item($pair[0], 6).Value = $pair[1]
}

# Insert the new row, carrying the formatting down from the row that will
# land at 42 (the old row 41) so the new row 41 keeps the same look (bold,
# bordered, centered index column etc.) as every other data row.
$ws4.Rows.Item(41).Insert()
$ws4.Range("A42").Copy()
$ws4.Range("A41").PasteSpecial(-4122)
$ws4.Application.CutCopyMode = $false

$ws4.Cells.Item(41, 1).Value = 40
$ws4.Cells.Item(41, 2).Value = "2024-08-24"
$ws4.Cells.Item(41, 3).Value = "广州·春日计划2024——特别二次元不插电音乐会"
$ws4.Cells.Item(41, 4).Value = "人民北路696号 广州友谊剧院"
$ws4.Cells.Item(41, 5).Value = "2024.08.24 19:30-08.24 21:00"
$ws4.Cells.Item(41, 6).Value = 12
$ws4.Cells.Item(41, 7).Value = 88
$ws4.Cells.Item(41, 8).Value = "https://show.bilibili.com/platform/detail.html?id=89964"
$ws4.Cells.Item(41, 9).Value = "//i0.hdslb.com/bfs/openplatform/202407/lHPV2n6t1722233858047.jpeg"
